# [CHG] INICIALIZANDO PÁGINA DE PLAN ESTRATEGICO
# Adds REGIÃO / SUPT / DICOM columns (P, Q, R) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header row (P1:R1) - bold font, centered/top aligned, thin left+right
# border (matches the existing header style used in A1:O1, but with a
# left/right-only border instead of a full box border).
# ---------------------------------------------------------------------
$p1 = $ws.Range("P1")
$p1.Value = "REGIÃO"
$p1.Font.Bold = $true
$p1.HorizontalAlignment = -4108
$p1.VerticalAlignment = -4160
$p1.Borders.Item(7).LineStyle = 1
$p1.Borders.Item(10).LineStyle = 1

# Re-use the freshly built header style for the other two header cells.
$ws.Range("P1").Copy()
$ws.Range("Q1:R1").PasteSpecial(-4122)
$ws.Range("Q1").Value = "SUPT"
$ws.Range("R1").Value = "DICOM"

# ---------------------------------------------------------------------
# Data rows 2-21: REGIÃO / SUPT / DICOM values
# ---------------------------------------------------------------------
$ws.Range("P2").Value = 13
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 6
$ws.Range("P3").Value = 6
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 1
$ws.Range("P4").Value = 18
$ws.Range("Q4").Value = 2
$ws.Range("R4").Value = 3
$ws.Range("P5").Value = 19
$ws.Range("Q5").Value = 4
$ws.Range("R5").Value = 2
$ws.Range("P6").Value = 12
$ws.Range("Q6").Value = 4
$ws.Range("R6").Value = 6
$ws.Range("P7").Value = 9
$ws.Range("Q7").Value = 4
$ws.Range("R7").Value = 7
$ws.Range("P8").Value = 6
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = 0
$ws.Range("P9").Value = 2
$ws.Range("Q9").Value = 3
$ws.Range("R9").Value = 4
$ws.Range("P10").Value = 14
$ws.Range("Q10").Value = 4
$ws.Range("R10").Value = 2
$ws.Range("P11").Value = 18
$ws.Range("Q11").Value = 0
$ws.Range("R11").Value = 5
$ws.Range("P12").Value = 9
$ws.Range("Q12").Value = 3
$ws.Range("R12").Value = 3
$ws.Range("P13").Value = 1
$ws.Range("Q13").Value = 2
$ws.Range("R13").Value = 2
$ws.Range("P14").Value = 4
$ws.Range("Q14").Value = 2
$ws.Range("R14").Value = 0
$ws.Range("P15").Value = 3
$ws.Range("Q15").Value = 0
$ws.Range("R15").Value = 4
$ws.Range("P16").Value = 14
$ws.Range("Q16").Value = 0
$ws.Range("R16").Value = 7
$ws.Range("P17").Value = 1
$ws.Range("Q17").Value = 3
$ws.Range("R17").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 1
$ws.Range("R18").Value = 6
$ws.Range("P19").Value = 13
$ws.Range("Q19").Value = 0
$ws.Range("R19").Value = 7
$ws.Range("P20").Value = 17
$ws.Range("Q20").Value = 1
$ws.Range("R20").Value = 0
$ws.Range("P21").Value = 7
$ws.Range("Q21").Value = 4
$ws.Range("R21").Value = 4

# ---------------------------------------------------------------------
# Reset the active selection back to A1 (default cursor position).
# ---------------------------------------------------------------------
$ws.Range("A1").Select()
